$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with changed values ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.543.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.186.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.186.60"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.92"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.511"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.710.61"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.540.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.186.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "513.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.42"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.90"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.07"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.06"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.53"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.86"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0894"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0422"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.15%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.858.19"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.85%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.38%  "

# --- Row reorders: Cosmos/Kaspa swap (rows 40-41), PEPE/dogwifhat swap (rows 43-44) ---
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.34%  "

$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.87"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.70%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0679"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.98%  "
